$d = $word.ActiveDocument

# 1. Remove the existing "_GoBack" bookmark that currently sits at the very
#    start of the document (it will be re-created at the new end of the
#    document content below).
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# 2. Append two new paragraphs after the current last paragraph:
#      - an empty paragraph (same style as the other spacer paragraphs,
#        sz/szCs = 36)
#      - a "Note : ..." paragraph (sz/szCs = 44) explaining that values
#        created with `const` cannot be reassigned.
#    The "_GoBack" bookmark is recreated at the very end of that last
#    paragraph, right after its final run.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$insertPoint = $d.Range($lastPara.Range.End, $lastPara.Range.End)

$newParagraphsXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/part" pkg:contentType="application/xml">
    <pkg:xmlData>
      <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:pPr>
          <w:jc w:val="center"/>
          <w:rPr>
            <w:sz w:val="36"/>
            <w:szCs w:val="36"/>
          </w:rPr>
        </w:pPr>
      </w:p>
      <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:pPr>
          <w:jc w:val="center"/>
          <w:rPr>
            <w:sz w:val="44"/>
            <w:szCs w:val="44"/>
          </w:rPr>
        </w:pPr>
        <w:proofErr w:type="gramStart"/>
        <w:r>
          <w:rPr>
            <w:sz w:val="44"/>
            <w:szCs w:val="44"/>
          </w:rPr>
          <w:t>Note :</w:t>
        </w:r>
        <w:proofErr w:type="gramEnd"/>
        <w:r>
          <w:rPr>
            <w:sz w:val="44"/>
            <w:szCs w:val="44"/>
          </w:rPr>
          <w:t xml:space="preserve"> Anything that changes the stored value of a variable created with the </w:t>
        </w:r>
        <w:proofErr w:type="spellStart"/>
        <w:r>
          <w:rPr>
            <w:sz w:val="44"/>
            <w:szCs w:val="44"/>
          </w:rPr>
          <w:t>const</w:t>
        </w:r>
        <w:proofErr w:type="spellEnd"/>
        <w:r>
          <w:rPr>
            <w:sz w:val="44"/>
            <w:szCs w:val="44"/>
          </w:rPr>
          <w:t xml:space="preserve"> keyword, Will result in an error. Like the Assignment/Variable Exclusive Mathematic operators.</w:t>
        </w:r>
        <w:bookmarkStart w:id="0" w:name="_GoBack"/>
        <w:bookmarkEnd w:id="0"/>
      </w:p>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

[void]$insertPoint.InsertXML($newParagraphsXml)
